$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.031.22"
$ws.Range("E2").Value = "  +6.87%  "
$ws.Range("D3").Value = "'1.742.78"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'228.52"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "'0.5441"
$ws.Range("E6").Value = "  +3.54%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").Value = "'0.06737"
$ws.Range("E9").Value = "  +5.91%  "
$ws.Range("D10").Value = "'21.73"
$ws.Range("E10").Value = "  +4.81%  "
$ws.Range("D11").Value = "'0.07783"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "'4.701"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "'1.741.49"
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("D14").Value = "'1.980.98"
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("D15").Value = "'0.5975"
$ws.Range("E15").Value = "  +5.91%  "
$ws.Range("D16").Value = "'0.0₅8384"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "'68.96"
$ws.Range("E17").Value = "  +5.25%  "
$ws.Range("D18").Value = "'28.011.19"
$ws.Range("E18").Value = "  +6.78%  "
$ws.Range("D19").Value = "'226.12"
$ws.Range("E19").Value = "  +17.61%  "
$ws.Range("D20").Value = "'4.843"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'10.91"
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("D23").Value = "'6.238"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'146.33"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "'0.1250"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").Value = "'7.478"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").Value = "'17.22"
$ws.Range("E28").Value = "  +7.85%  "
$ws.Range("D29").Value = "'1.647"
$ws.Range("E29").Value = "  +10.02%  "
$ws.Range("D30").Value = "'0.05676"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").Value = "'3.706"
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("D33").Value = "'3.521"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "'1.679"
$ws.Range("E34").Value = "  +5.96%  "
$ws.Range("D35").Value = "'0.9849"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").Value = "'2.864"
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").Value = "'2.450"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").Value = "'0.5962"
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").Value = "'0.01671"
$ws.Range("E39").Value = "  +4.46%  "
$ws.Range("D40").Value = "'5.958"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "'1.051.67"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("D42").Value = "'0.8492"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'102.10"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'1.886.70"
$ws.Range("E45").Value = "  +5.05%  "
$ws.Range("E46").Value = "  +4.65%  "
$ws.Range("D47").Value = "'60.03"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").Value = "'8.290"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'0.05324"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("E51").Value = "  -0.81%  "
